$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "31.245.37"
$ws.Range("E2").Value = "  +4.53%  "

# Row 3
$ws.Range("D3").Value = "1.701.50"
$ws.Range("E3").Value = "  +4.02%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.45"
$ws.Range("E5").Value = "  +2.80%  "

# Row 6
$ws.Range("E6").Value = "  +2.88%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.76"
$ws.Range("E8").Value = "  +3.29%  "

# Row 9
$ws.Range("E9").Value = "  +2.86%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0643"
$ws.Range("E10").Value = "  +5.59%  "

# Row 11
$ws.Range("E11").Value = "  +1.25%  "

# Row 12
$ws.Range("D12").Value = "1.943.58"
$ws.Range("E12").Value = "  +3.88%  "

# Row 13
$ws.Range("D13").Value = "1.704.33"
$ws.Range("E13").Value = "  +4.15%  "

# Row 14
$ws.Range("E14").Value = "  +3.75%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.14"
$ws.Range("E15").Value = "  +7.37%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.17"
$ws.Range("E16").Value = "  +8.15%  "

# Row 17
$ws.Range("D17").Value = "31.233.16"
$ws.Range("E17").Value = "  +4.42%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.10"
$ws.Range("E18").Value = "  +3.80%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "250.11"
$ws.Range("E19").Value = "  +3.89%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0725"
$ws.Range("E20").Value = "  +3.14%  "

# Row 21
$ws.Range("E21").Value = "  +0.11%  "

# Row 22
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.18"
$ws.Range("E22").Value = "  +2.88%  "

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.28"
$ws.Range("E23").Value = "  +3.28%  "

# Row 24
$ws.Range("E24").Value = "  -0.23%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.73"
$ws.Range("E25").Value = "  +0.93%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.02"
$ws.Range("E26").Value = "  +3.09%  "

# Row 27
$ws.Range("E27").Value = "  +3.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.77"
$ws.Range("E28").Value = "  +2.08%  "

# Row 29
$ws.Range("E29").Value = "  -0.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.78"
$ws.Range("E30").Value = "  +11.58%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0504"
$ws.Range("E31").Value = "  +2.56%  "

# Row 32
$ws.Range("E32").Value = "  +3.62%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.40"
$ws.Range("E33").Value = "  +6.37%  "

# Row 34
$ws.Range("D34").Value = "1.511.41"
$ws.Range("E34").Value = "  +6.00%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.74"
$ws.Range("E35").Value = "  +2.93%  "

# Row 36
$ws.Range("E36").Value = "  +1.19%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.617"
$ws.Range("E37").Value = "  +9.98%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "83.19"
$ws.Range("E38").Value = "  +9.13%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0181"
$ws.Range("E39").Value = "  +4.79%  "

# Row 40
$ws.Range("E40").Value = "  -3.39%  "

# Row 41
$ws.Range("E41").Value = "  +0.45%  "

# Row 42
$ws.Range("E42").Value = "  +3.51%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.853"
$ws.Range("E43").Value = "  +2.37%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0503"
$ws.Range("E44").Value = "  +0.71%  "

# Row 45
$ws.Range("E45").Value = "  +2.96%  "

# Row 46
$ws.Range("E46").Value = "  +0.00%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.33"
$ws.Range("E47").Value = "  +7.10%  "

# Row 48
$ws.Range("E48").Value = "  +3.35%  "

# Row 49
$ws.Range("D49").Value = "1.834.15"
$ws.Range("E49").Value = "  +3.16%  "

# Row 50
$ws.Range("E50").Value = "  +7.73%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "94.14"
$ws.Range("E51").Value = "  +1.23%  "
